$wb = $excel.ActiveWorkbook

$wsAceptadas = $wb.Worksheets.Item("VENTAS CONSULTA ACEPTADAS")
$wsRechazadas = $wb.Worksheets.Item("VENTAS CONSULTAS RECHAZADAS")

# --- Sheet "VENTAS CONSULTA ACEPTADAS" ---
# Duplicate the data row (row 3) into row 4, carrying values + formatting.
$wsAceptadas.Range("A3:AH3").Copy($wsAceptadas.Range("A4:AH4"))

# Taller rows to fit the wrapped text now shown in both rows.
$wsAceptadas.Rows(3).RowHeight = 45.5
$wsAceptadas.Rows(4).RowHeight = 45.5

# New sale reference number for the duplicated row.
$wsAceptadas.Range("O4").Value = "5550455"

# --- Sheet "VENTAS CONSULTAS RECHAZADAS" ---
# Zoom this sheet in before returning focus to the main sheet.
$wsRechazadas.Activate()
$excel.ActiveWindow.Zoom = 115
$wsRechazadas.Range("D15").Select()

# --- back to "VENTAS CONSULTA ACEPTADAS" ---
$wsAceptadas.Activate()
$excel.ActiveWindow.Zoom = 160
$wsAceptadas.Range("O4").Select()
